$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial value that needs to move
# forward by one day (45179 -> 45180) for every data row (rows 2-358).
$ws.Range("C2:C358").Value = 45180
